$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ticu1245", "Y", "a", "Y"),
    @("orej1242", "N", "b", "N"),
    @("nade1244", "N", "c", "Y"),
    @("mara1409", "N", "a", "N")
)

$row = 4
foreach ($rowData in $data) {
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $row++
}
